$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 401757.4
$ws.Range("J17").Value = 401757.4
$ws.Range("L17").Value = 1205272.2
$ws.Range("N17").Value = -1205608.2

$ws.Range("H39").Value = 4183.1
$ws.Range("I39").Value = 226.4
$ws.Range("J39").Value = 8139.8
$ws.Range("K39").Value = 679.2
$ws.Range("L39").Value = 24419.4
$ws.Range("M39").Value = -383.2
$ws.Range("N39").Value = -25011.4

$ws.Range("H62").Value = 5381.2
$ws.Range("J62").Value = 5953
$ws.Range("L62").Value = 5953
$ws.Range("N62").Value = -7201

$ws.Range("H65").Value = 5381.2
$ws.Range("J65").Value = 5953
$ws.Range("L65").Value = 29765
$ws.Range("N65").Value = -36005

$ws.Range("H70").Value = 7424.9375
$ws.Range("J70").Value = 21750
$ws.Range("L70").Value = 65250
$ws.Range("N70").Value = -65790

$ws.Range("H73").Value = 7424.9375
$ws.Range("J73").Value = 21750
$ws.Range("L73").Value = 65250
$ws.Range("N73").Value = -67122

$ws.Range("H96").Value = 1785
$ws.Range("I96").Value = 1198.4
$ws.Range("J96").Value = 3251.5
$ws.Range("K96").Value = 3595.2
$ws.Range("L96").Value = 9754.5
$ws.Range("M96").Value = -2222.2
$ws.Range("N96").Value = -12500.5

$ws.Range("H111").Value = 744
$ws.Range("I111").Value = 200
$ws.Range("J111").Value = 1016
$ws.Range("K111").Value = 600
$ws.Range("L111").Value = 3048
$ws.Range("M111").Value = 2467
$ws.Range("N111").Value = -9182

$ws.Range("H113").Value = 3351.2
$ws.Range("I113").Value = 2383.3333
$ws.Range("J113").Value = 4803
$ws.Range("K113").Value = 2383.3333
$ws.Range("L113").Value = 4803
$ws.Range("M113").Value = 870.6667000000002
$ws.Range("N113").Value = -11311

$ws.Range("H137").Value = 1674.079
$ws.Range("I137").Value = 1145.2759
$ws.Range("J137").Value = 3378
$ws.Range("K137").Value = 3435.8277
$ws.Range("L137").Value = 10134
$ws.Range("M137").Value = -885.8277000000003
$ws.Range("N137").Value = -15234

$ws.Range("H138").Value = 3510.198
$ws.Range("J138").Value = 3622.859
$ws.Range("L138").Value = 10868.577
$ws.Range("N138").Value = -21148.577

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 466346.53
$ws.Range("I2").Value = 1113315.8
$ws.Range("J2").Value = 4225.619
$ws.Range("K2").Value = 1113315.8
$ws.Range("L2").Value = 4225.619
$ws.Range("M2").Value = -1113202.8
$ws.Range("N2").Value = -4451.619

$ws.Range("H32").Value = 11317.5
$ws.Range("I32").Value = 6813.543
$ws.Range("J32").Value = 18171.348
$ws.Range("K32").Value = 6813.543
$ws.Range("L32").Value = 18171.348
$ws.Range("M32").Value = -6526.543
$ws.Range("N32").Value = -18745.348

$ws.Range("H45").Value = 1151.8462
$ws.Range("I45").Value = 979.1
$ws.Range("K45").Value = 979.1
$ws.Range("M45").Value = -602.1

$ws.Range("H97").Value = 530.1818
$ws.Range("I97").Value = 558.2
$ws.Range("K97").Value = 558.2
$ws.Range("M97").Value = -62.20000000000005

$ws.Range("H101").Value = 199854.62
$ws.Range("J101").Value = 199854.62
$ws.Range("L101").Value = 199854.62
$ws.Range("N101").Value = -206344.62

$ws.Range("H104").Value = 1648369
$ws.Range("J104").Value = 1648369
$ws.Range("L104").Value = 1648369
$ws.Range("N104").Value = -1655357

$ws.Range("H116").Value = 466346.53
$ws.Range("I116").Value = 1113315.8
$ws.Range("J116").Value = 4225.619
$ws.Range("K116").Value = 1113315.8
$ws.Range("L116").Value = 4225.619
$ws.Range("M116").Value = -1111021.8
$ws.Range("N116").Value = -8813.618999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 466346.53
$ws.Range("I3").Value = 1113315.8
$ws.Range("J3").Value = 4225.619
$ws.Range("K3").Value = 1113315.8
$ws.Range("L3").Value = 4225.619
$ws.Range("M3").Value = -1113201.8
$ws.Range("N3").Value = -4453.619

$ws.Range("H99").Value = 2213.875
$ws.Range("J99").Value = 2702.2
$ws.Range("L99").Value = 2702.2
$ws.Range("N99").Value = -5698.2

$ws.Range("H105").Value = 5200.727
$ws.Range("I105").Value = 5230.9
$ws.Range("K105").Value = 5230.9
$ws.Range("M105").Value = -3483.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 11218.375
$ws.Range("J95").Value = 11218.375
$ws.Range("L95").Value = 11218.375
$ws.Range("N95").Value = -16710.375

$ws.Range("H96").Value = 5748.5
$ws.Range("J96").Value = 5748.5
$ws.Range("L96").Value = 5748.5
$ws.Range("N96").Value = -11240.5

$ws.Range("H122").Value = 4134.7544
$ws.Range("I122").Value = 3993.2126
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 11979.6378
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -9529.6378
$ws.Range("N122").Value = -19300

$ws.Range("H132").Value = 500001500
$ws.Range("I132").Value = 500001500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1500004500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1500001970
$ws.Range("N132").ClearContents()

$ws.Range("H141").Value = 278122.84
$ws.Range("I141").Value = 57140
$ws.Range("J141").Value = 416237.12
$ws.Range("K141").Value = 57140
$ws.Range("L141").Value = 416237.12
$ws.Range("M141").Value = -51960
$ws.Range("N141").Value = -426597.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 162.11765
$ws.Range("I12").Value = 27.666666
$ws.Range("J12").Value = 235.45454
$ws.Range("K12").Value = 82.99999800000001
$ws.Range("L12").Value = 706.3636200000001
$ws.Range("M12").Value = 90.00000199999999
$ws.Range("N12").Value = -1052.36362

$ws.Range("H107").Value = 1640.2307
$ws.Range("I107").Value = 144.5
$ws.Range("J107").Value = 1912.1818
$ws.Range("K107").Value = 433.5
$ws.Range("L107").Value = 5736.5454
$ws.Range("M107").Value = 1486.5
$ws.Range("N107").Value = -9576.545399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2741.3438
$ws.Range("I80").Value = 3121.5
$ws.Range("J80").Value = 2653.6155
$ws.Range("K80").Value = 3121.5
$ws.Range("L80").Value = 2653.6155
$ws.Range("M80").Value = -2123.5
$ws.Range("N80").Value = -4649.6155

$ws.Range("H83").Value = 2741.3438
$ws.Range("I83").Value = 3121.5
$ws.Range("J83").Value = 2653.6155
$ws.Range("K83").Value = 15607.5
$ws.Range("L83").Value = 13268.0775
$ws.Range("M83").Value = -10615.5
$ws.Range("N83").Value = -23252.0775

$ws.Range("H98").Value = 25285.572
$ws.Range("J98").Value = 25285.572
$ws.Range("L98").Value = 25285.572
$ws.Range("N98").Value = -31275.572

$ws.Range("H99").Value = 19445.777
$ws.Range("I99").Value = 7859
$ws.Range("K99").Value = 7859
$ws.Range("M99").Value = -5613

$ws.Range("H122").Value = 96627.46000000001
$ws.Range("I122").Value = 111514.27
$ws.Range("K122").Value = 334542.81
$ws.Range("M122").Value = -332092.81

$ws.Range("H126").Value = 9408.909
$ws.Range("I126").Value = 7566.2856
$ws.Range("J126").Value = 12633.5
$ws.Range("K126").Value = 22698.8568
$ws.Range("L126").Value = 37900.5
$ws.Range("M126").Value = -20228.8568
$ws.Range("N126").Value = -42840.5

$ws.Range("H132").Value = 2979618.5
$ws.Range("I132").Value = 3791372.5
$ws.Range("K132").Value = 11374117.5
$ws.Range("M132").Value = -11371587.5

$ws.Range("H135").Value = 133473.75
$ws.Range("J135").Value = 133473.75
$ws.Range("L135").Value = 133473.75
$ws.Range("N135").Value = -143613.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2420.5386
$ws.Range("I16").Value = 738.2857
$ws.Range("K16").Value = 738.2857
$ws.Range("M16").Value = -568.2857

$ws.Range("H99").Value = 38253.145
$ws.Range("I99").Value = 25897.6
$ws.Range("K99").Value = 25897.6
$ws.Range("M99").Value = -22902.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23659.334
$ws.Range("J54").Value = 23642.215
$ws.Range("L54").Value = 23642.215
$ws.Range("N54").Value = -24682.215

$ws.Range("H105").Value = 18857
$ws.Range("J105").Value = 18857
$ws.Range("L105").Value = 18857
$ws.Range("N105").Value = -25845

$ws.Range("H110").Value = 258999
$ws.Range("J110").Value = 258999
$ws.Range("L110").Value = 258999
$ws.Range("N110").Value = -267179

$ws.Range("H132").Value = 20008218
$ws.Range("I132").Value = 33338664
$ws.Range("K132").Value = 100015992
$ws.Range("M132").Value = -100013462

$ws.Range("H136").Value = 20835678
$ws.Range("I136").Value = 26317180
$ws.Range("K136").Value = 78951540
$ws.Range("M136").Value = -78948990

$ws.Range("H137").Value = 125000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 125000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 125000
$ws.Range("N137").Value = -135200
$ws.Range("M137").ClearContents()
